$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("BQ1").Value2 = "LocMinDed6All"
$ws.Range("BQ2").Formula = "=0"
Write-Host "Before replace: $($ws.Range('BQ2').Text)"
$rng = $ws.Range("BQ2")
$rng.Replace("=0", "=IFNA(HLOOKUP(BQ`$1,[1]location!`$A`$1:`$X`$100,`$B2,FALSE),0)", 1, 1, $false, $false, $false)
Write-Host "After replace formula: $($ws.Range('BQ2').Formula)"
Write-Host "After replace value: $($ws.Range('BQ2').Text)"
